$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.990.72"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "1.559.59"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("E4").Value = "  +0.19%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.02"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.62%  "

$ws.Range("E6").Value = "  +0.62%  "

$ws.Range("E7").Value = "  +0.21%  "

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.06"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  +0.79%  "

$ws.Range("E10").Value = "  +1.90%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0854"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  -0.26%  "

$ws.Range("D12").Value = "1.781.88"
$ws.Range("E12").Value = "  +0.59%  "

$ws.Range("D13").Value = "1.557.65"
$ws.Range("E13").Value = "  +0.65%  "

$ws.Range("E14").Value = "  -0.30%  "

$ws.Range("E15").Value = "  +0.16%  "

$ws.Range("D16").Value = "27.009.31"
$ws.Range("E16").Value = "  +0.31%  "

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.78"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("E18").Value = "  +1.53%  "

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.45"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -0.90%  "

$ws.Range("E20").Value = "  +1.48%  "

$ws.Range("E21").Value = "  +0.17%  "

$ws.Range("E22").Value = "  +2.29%  "

$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("E24").Value = "  -0.73%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.77"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -0.87%  "

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.60"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -0.23%  "

$ws.Range("E27").Value = "  +0.78%  "

$ws.Range("E28").Value = "  +1.60%  "

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("E30").Value = "  +1.33%  "

$ws.Range("E31").Value = "  +3.35%  "

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.18"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +3.81%  "

$ws.Range("D34").Value = "1.423.41"
$ws.Range("E34").Value = "  +0.17%  "

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.08"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +11.05%  "

$ws.Range("E36").Value = "  +1.12%  "

$ws.Range("E37").Value = "  +2.48%  "

$ws.Range("E38").Value = "  +1.05%  "

$ws.Range("E39").Value = "  +1.84%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.809"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +0.28%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.78"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +0.08%  "

$ws.Range("E42").Value = "  +0.25%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +0.70%  "

$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.31"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("E45").Value = "  +0.49%  "

$ws.Range("E46").Value = "  -1.12%  "

$ws.Range("D47").Value = "1.696.13"
$ws.Range("E47").Value = "  +0.63%  "

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.74"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -1.06%  "

$ws.Range("E49").Value = "  +2.94%  "

$ws.Range("E50").Value = "  -0.32%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0960"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +1.08%  "
